# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Timestamp footer text (A1)
$ws.Range("A1").Value = "Datos actualizados a 31 de Agosto de 2020 a las 17:26"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 6180197
$ws.Range("C4").Value = 6961
$ws.Range("D4").Value = 3425938
$ws.Range("E4").Value = 2566894
$ws.Range("G4").Value = 141
$ws.Range("H4").Value = 187365

# India (row 6)
$ws.Range("B6").Value = 3649639
$ws.Range("C6").Value = 30470
$ws.Range("D6").Value = 2800671
$ws.Range("E6").Value = 784017
$ws.Range("G6").Value = 334
$ws.Range("H6").Value = 64951

# Reino Unido (row 16)
$ws.Range("B16").Value = 335873
$ws.Range("C16").Value = 1406
$ws.Range("G16").Value = 2
$ws.Range("H16").Value = 41501

# Kenia (row 68)
$ws.Range("B68").Value = 34201
$ws.Range("C68").Value = 144
$ws.Range("D68").Value = 19893
$ws.Range("E68").Value = 13731
$ws.Range("G68").Value = 3
$ws.Range("H68").Value = 577

# Sri Lanka (row 127)
$ws.Range("B127").Value = 3018
$ws.Range("C127").Value = 6
$ws.Range("E127").Value = 138

# Jordania overtakes Sierra Leona and Aruba in the ranking (rows 142-144
# keep their position in the table but the countries attached to them
# shift down, carrying each other's previous figures along).
$ws.Range("A142").Value = "Jordania"
$ws.Range("B142").Value = 2034
$ws.Range("C142").Value = 68
$ws.Range("D142").Value = 1508
$ws.Range("E142").Value = 511
$ws.Range("H142").Value = 15

$ws.Range("A143").Value = "Sierra Leona"
$ws.Range("B143").Value = 2022
$ws.Range("D143").Value = 1594
$ws.Range("E143").Value = 358
$ws.Range("H143").Value = 70

$ws.Range("A144").Value = "Aruba"
$ws.Range("B144").Value = 1997
$ws.Range("D144").Value = 765
$ws.Range("E144").Value = 1222
$ws.Range("H144").Value = 10

# Burkina Faso (row 156)
$ws.Range("B156").Value = 1368
$ws.Range("C156").Value = 11
$ws.Range("E156").Value = 255
